$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.651.48'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.244.16'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').Value = "'306.63"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = "'94.29"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D9').Value = "'0.515"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('D10').Value = "'34.79"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '2.586.80'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').Value = '2.243.29'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = "'0.830"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = "'13.55"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').Value = '44.421.21'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('D20').Value = "'11.78"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('E21').Value = '  -3.20%  '
$ws.Range('D22').Value = "'65.30"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').Value = "'237.74"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  +3.73%  '
$ws.Range('D28').Value = "'9.77"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('D29').Value = "'36.92"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.55%  '
$ws.Range('D30').Value = "'19.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').Value = "'148.53"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('D33').Value = "'0.0781"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').Value = "'3.18"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('E38').Value = '  +5.34%  '
$ws.Range('D39').Value = "'15.13"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.60%  '
$ws.Range('D40').Value = "'3.34"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('D42').Value = "'0.0298"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '1.810.14'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('D45').Value = "'1.79"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.76%  '
$ws.Range('D46').Value = "'81.87"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('E47').Value = '  -2.00%  '
$ws.Range('D48').Value = "'97.90"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = "'69.13"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'4.82"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').Value = "'53.94"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
